$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '[Gunther%Eysenbach%NULL%0,               Nikki%Lam%NULL%2,               Nikki%Lam%NULL%0,               Benjamin%Woo%NULL%1,               Kylie%Zeng%kzeng2@student.touro.edu%2,               Kylie%Zeng%kzeng2@student.touro.edu%0,               Stephanie N%Bernardo%NULL%2,               Stephanie N%Bernardo%NULL%0,               Weldon E%Havins%NULL%2,               Weldon E%Havins%NULL%0]'
$ws.Range('E3').Value = '[Gunther%Eysenbach%NULL%0,               Guy%Fagherazzi%NULL%1,               Juan%Yang%NULL%0,               Juan%Yang%NULL%0,               Claire%Wang%NULL%1,               Céline%Mavrot%NULL%1,               Carlos K H%Wong%NULL%2,               Carlos K H%Wong%NULL%0,               Janet Y H%Wong%NULL%2,               Janet Y H%Wong%NULL%0,               Eric H M%Tang%NULL%2,               Eric H M%Tang%NULL%0,               Chi Ho%Au%NULL%2,               Chi Ho%Au%NULL%0,               Kristy T K%Lau%NULL%2,               Kristy T K%Lau%NULL%0,               Abraham K C%Wai%awai@hku.hk%2,               Abraham K C%Wai%awai@hku.hk%0]'
$ws.Range('E4').Value = '[Mark J.%Siedner%NULL%1,               Guy%Harling%NULL%2,               Guy%Harling%NULL%0,               Zahra%Reynolds%NULL%2,               Zahra%Reynolds%NULL%0,               Rebecca F.%Gilbert%NULL%1,               Sebastien%Haneuse%NULL%2,               Sebastien%Haneuse%NULL%0,               Atheendar S.%Venkataramani%NULL%2,               Atheendar S.%Venkataramani%NULL%0,               Alexander C.%Tsai%NULL%2,               Alexander C.%Tsai%NULL%0,               Donald A.%Redelmeier%NULL%3,               Donald A.%Redelmeier%NULL%0,               Donald A.%Redelmeier%NULL%0,               NULL%NULL%NULL%0,               NULL%NULL%NULL%0,               NULL%NULL%NULL%0]'
$ws.Range('E5').Value = '[ B.%Krishnamachari%null%1,                A.% Dsida%null%1,                D.% Zastrow%null%1,                B.% Harper%null%1,                A.% Morris%null%1,                A. % Santella%null%1]'
$ws.Range('E6').Value = '[ K. P.%Kepp%null%1,                C. % Bjornskov%null%1]'
$ws.Range('E7').Value = '[An%Pan%xref no email%1,        Li%Liu%xref no email%1,        Chaolong%Wang%xref no email%1,        Huan%Guo%xref no email%1,        Xingjie%Hao%xref no email%1,        Qi%Wang%xref no email%1,        Jiao%Huang%xref no email%1,        Na%He%xref no email%1,        Hongjie%Yu%xref no email%1,        Xihong%Lin%xref no email%1,        Sheng%Wei%xref no email%1,        Tangchun%Wu%xref no email%1]'
$ws.Range('E8').Value = '[Ke-wei%Wang%NULL%1,               Jie%Gao%NULL%1,               Hua%Wang%NULL%1,               Xiao-long%Wu%NULL%1,               Qin-fang%Yuan%NULL%1,               Fei-yu%Guo%NULL%1,               Zhi-jie%Zhang%NULL%1,               Yang%Cheng%woerseng@126.com%1]'
$ws.Range('E9').Value = '[Benjamin J%Cowling%NULL%0,               Sheikh Taslim%Ali%NULL%0,               Tiffany W Y%Ng%NULL%0,               Tim K%Tsang%NULL%0,               Julian C M%Li%NULL%0,               Min Whui%Fong%NULL%0,               Qiuyan%Liao%NULL%0,               Mike YW%Kwan%NULL%0,               So Lun%Lee%NULL%0,               Susan S%Chiu%NULL%0,               Joseph T%Wu%NULL%0,               Peng%Wu%pengwu@hku.hk%0,               Gabriel M%Leung%NULL%0]'
$ws.Range('E10').Value = '[Vincent Chi-Chung%Cheng%NULL%1,               Shuk-Ching%Wong%NULL%1,               Vivien Wai-Man%Chuang%NULL%1,               Simon Yung-Chun%So%NULL%1,               Jonathan Hon-Kwan%Chen%NULL%1,               Siddharth%Sridhar%NULL%1,               Kelvin Kai-Wang%To%NULL%0,               Jasper Fuk-Woo%Chan%NULL%0,               Ivan Fan-Ngai%Hung%NULL%0,               Pak-Leung%Ho%NULL%1,               Kwok-Yung%Yuen%NULL%0]'
$ws.Range('E11').Value = '[Tran Phuoc Bao%Thu%NULL%1,               Pham Nguyen Hong%Ngoc%NULL%1,               Nguyen Minh%Hai%NULL%1,               Le Anh%Tuan%NULL%1]'
$ws.Range('E12').Value = '[Hye%Seong%NULL%1,               Hak Jun%Hyun%NULL%1,               Jin Gu%Yun%NULL%1,               Ji Yun%Noh%NULL%1,               Hee Jin%Cheong%NULL%0,               Woo Joo%Kim%NULL%0,               Joon Young%Song%NULL%0]'
$ws.Range('E13').Value = '[Sultan Ayoub%Meo%sultanmeo@hotmail.com%1,               Abdulelah Adnan%Abukhalaf%Abdulelahabukhalaf@gmail.com%1,               Ali Abdullah%Alomar%AliAlomarMD@gmail.com%1,               Faris Jamal%AlMutairi%faris11300@gmail.com%1,               Adnan Mahmood%Usmani%adnan_mahmood71@hotmail.com%1,               David C.%Klonoff%dklonoff@diabetestechnology.org%1]'
